$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 6235
$ws.Range("C21").Value = 986
$ws.Range("D21").Value = 5608035
$ws.Range("E21").Value = 899.4442662389736
$ws.Range("F21").Value = 8.227738239888915
$ws.Range("G21").Value = 4.008438818565407
$ws.Range("H21").Value = 28.00892130787414
